$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Vtn"
$ws.Range("C2").Value = "Itgb6"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 7.134618
$ws.Range("H2").Value = 21.403854
$ws.Range("I2").Value = 0.0965317920926077
$ws.Range("J2").Value = 0.0965317920926077
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.069456
$ws.Range("N2").Value = 0.208368
$ws.Range("O2").Value = 0.01627409789654661
$ws.Range("P2").Value = 0.01627409789654661
$ws.Range("Q2").Value = 0.495542027808
$ws.Range("R2").Value = 4.459878250271999
$ws.Range("S2").Value = 0.001570967834644181
$ws.Range("T2").Value = 0.001570967834644181

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Vtn"
$ws.Range("C3").Value = "Itgb6"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 7.134618
$ws.Range("H3").Value = 21.403854
$ws.Range("I3").Value = 0.0965317920926077
$ws.Range("J3").Value = 0.0965317920926077
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.6957970000000001
$ws.Range("N3").Value = 2.087391
$ws.Range("O3").Value = 0.1630308179872645
$ws.Range("P3").Value = 0.1630308179872644
$ws.Range("Q3").Value = 4.964245800546
$ws.Range("R3").Value = 44.678212204914
$ws.Range("S3").Value = 0.01573765702663438
$ws.Range("T3").Value = 0.01573765702663438

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Vtn"
$ws.Range("C4").Value = "Itgb6"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 7.134618
$ws.Range("H4").Value = 21.403854
$ws.Range("I4").Value = 0.0965317920926077
$ws.Range("J4").Value = 0.0965317920926077
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.502633333333333
$ws.Range("N4").Value = 10.5079
$ws.Range("O4").Value = 0.820695084116189
$ws.Range("P4").Value = 0.820695084116189
$ws.Range("Q4").Value = 24.9899508274
$ws.Range("R4").Value = 224.9095574466
$ws.Range("S4").Value = 0.07922316723132915
$ws.Range("T4").Value = 0.07922316723132915

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Vtn"
$ws.Range("C5").Value = "Itgb6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 17.50798033333334
$ws.Range("H5").Value = 52.52394100000001
$ws.Range("I5").Value = 0.2368839813846793
$ws.Range("J5").Value = 0.2368839813846794
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.069456
$ws.Range("N5").Value = 0.208368
$ws.Range("O5").Value = 0.01627409789654661
$ws.Range("P5").Value = 0.01627409789654661
$ws.Range("Q5").Value = 1.216034282032
$ws.Range("R5").Value = 10.944308538288
$ws.Range("S5").Value = 0.003855073103177995
$ws.Range("T5").Value = 0.003855073103177996

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Vtn"
$ws.Range("C6").Value = "Itgb6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 17.50798033333334
$ws.Range("H6").Value = 52.52394100000001
$ws.Range("I6").Value = 0.2368839813846793
$ws.Range("J6").Value = 0.2368839813846794
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.6957970000000001
$ws.Range("N6").Value = 2.087391
$ws.Range("O6").Value = 0.1630308179872645
$ws.Range("P6").Value = 0.1630308179872644
$ws.Range("Q6").Value = 12.18200019199234
$ws.Range("R6").Value = 109.638001727931
$ws.Range("S6").Value = 0.0386193892532242
$ws.Range("T6").Value = 0.0386193892532242

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Vtn"
$ws.Range("C7").Value = "Itgb6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 17.50798033333334
$ws.Range("H7").Value = 52.52394100000001
$ws.Range("I7").Value = 0.2368839813846793
$ws.Range("J7").Value = 0.2368839813846794
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.502633333333333
$ws.Range("N7").Value = 10.5079
$ws.Range("O7").Value = 0.820695084116189
$ws.Range("P7").Value = 0.820695084116189
$ws.Range("Q7").Value = 61.32403551487779
$ws.Range("R7").Value = 551.9163196339
$ws.Range("S7").Value = 0.1944095190282772
$ws.Range("T7").Value = 0.1944095190282772

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Vtn"
$ws.Range("C8").Value = "Itgb6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 49.26691733333334
$ws.Range("H8").Value = 147.800752
$ws.Range("I8").Value = 0.6665842265227129
$ws.Range("J8").Value = 0.666584226522713
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.069456
$ws.Range("N8").Value = 0.208368
$ws.Range("O8").Value = 0.01627409789654661
$ws.Range("P8").Value = 0.01627409789654661
$ws.Range("Q8").Value = 3.421883010304001
$ws.Range("R8").Value = 30.796947092736
$ws.Range("S8").Value = 0.01084805695872443
$ws.Range("T8").Value = 0.01084805695872443

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Vtn"
$ws.Range("C9").Value = "Itgb6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 49.26691733333334
$ws.Range("H9").Value = 147.800752
$ws.Range("I9").Value = 0.6665842265227129
$ws.Range("J9").Value = 0.666584226522713
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.6957970000000001
$ws.Range("N9").Value = 2.087391
$ws.Range("O9").Value = 0.1630308179872645
$ws.Range("P9").Value = 0.1630308179872644
$ws.Range("Q9").Value = 34.27977327978135
$ws.Range("R9").Value = 308.5179595180321
$ws.Range("S9").Value = 0.1086737717074059
$ws.Range("T9").Value = 0.1086737717074059

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Vtn"
$ws.Range("C10").Value = "Itgb6"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 49.26691733333334
$ws.Range("H10").Value = 147.800752
$ws.Range("I10").Value = 0.6665842265227129
$ws.Range("J10").Value = 0.666584226522713
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.502633333333333
$ws.Range("N10").Value = 10.5079
$ws.Range("O10").Value = 0.820695084116189
$ws.Range("P10").Value = 0.820695084116189
$ws.Range("Q10").Value = 172.5639468823111
$ws.Range("R10").Value = 1553.0755219408
$ws.Range("S10").Value = 0.5470623978565826
$ws.Range("T10").Value = 0.5470623978565827
